$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "0.99915139503728889"
$ws.Range("B1").Value = "0.00084860496271118668"
$ws.Range("C1").Value = "0.99915139503728889"
$ws.Range("D1").Value = "0.00084860496271118733"
$ws.Range("E1").Value = "0.99915139503728867"
$ws.Range("F1").Value = "0.00084860496271118668"
$ws.Range("G1").Value = "0.00084860496271118668"
$ws.Range("H1").Value = "0.99915139503728889"
$ws.Range("I1").Value = "0.99915139503728889"
$ws.Range("J1").Value = "0.00084860496271126952"
$ws.Range("K1").Value = "0.00084860496271121683"
$ws.Range("L1").Value = "0.99915139503728889"
$ws.Range("M1").Value = "0.99915139503728889"
$ws.Range("N1").Value = "0.99915139503728889"
$ws.Range("O1").Value = "0.00084860496271114147"
$ws.Range("P1").Value = "0.00084860496271121227"
$ws.Range("Q1").Value = "0.00084860496271118733"
$ws.Range("R1").Value = "0.00084860496271118668"
$ws.Range("S1").Value = "0.99915139503728889"
$ws.Range("T1").Value = "0.99915139503728889"

$ws.Range("A2").Value = "0.00084860496271118668"
$ws.Range("B2").Value = "0.99915139503728889"
$ws.Range("C2").Value = "0.00084860496271118733"
$ws.Range("D2").Value = "0.99915139503728889"
$ws.Range("E2").Value = "0.00084860496271114668"
$ws.Range("F2").Value = "0.99915139503728889"
$ws.Range("G2").Value = "0.99915139503728889"
$ws.Range("H2").Value = "0.00084860496271118733"
$ws.Range("I2").Value = "0.00084860496271118733"
$ws.Range("J2").Value = "0.99915139503728889"
$ws.Range("K2").Value = "0.99915139503728889"
$ws.Range("L2").Value = "0.00084860496271119417"
$ws.Range("M2").Value = "0.00084860496271118668"
$ws.Range("N2").Value = "0.00084860496271118668"
$ws.Range("O2").Value = "0.99915139503728867"
$ws.Range("P2").Value = "0.99915139503728889"
$ws.Range("Q2").Value = "0.99915139503728889"
$ws.Range("R2").Value = "0.99915139503728889"
$ws.Range("S2").Value = "0.00084860496271128817"
$ws.Range("T2").Value = "0.00084860496271118668"

$ws.Range("A3").Value = "0.94512698668994977"
$ws.Range("B3").Value = "0.99953755535884092"
$ws.Range("C3").Value = "0.026030005273385688"
$ws.Range("D3").Value = "0.00013587496717803083"
$ws.Range("E3").Value = "0.047584009735105487"
$ws.Range("F3").Value = "0.000023218843091358534"

$ws.Range("A5").Value = "5"
$ws.Range("B5").Value = "1"
$ws.Range("C5").Value = "0"
$ws.Range("D5").Value = "1"
$ws.Range("E5").Value = "0"
$ws.Range("F5").Value = "5"

$ws.Range("A6").Value = "1"
$ws.Range("B6").Value = "3"
$ws.Range("C6").Value = "2"
$ws.Range("D6").Value = "1"
$ws.Range("E6").Value = "2"
$ws.Range("F6").Value = "2"

$ws.Range("B9").Value = "5"

$ws.Range("B10").Value = "1"

$ws.Range("B11").Value = "0"

$ws.Range("B12").Value = "1"

$ws.Range("B13").Value = "0"

$ws.Range("B14").Value = "5"

$ws.Range("A1:J2").Select()
